$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

Set-TextCell 2 2 'Bitcoin'
Set-TextCell 2 3 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextCell 2 4 '28.998.42'
Set-TextCell 2 5 '  -1.81%  '
Set-TextCell 3 2 'Ethereum'
Set-TextCell 3 3 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextCell 3 4 '1.908.26'
Set-TextCell 3 5 '  -3.39%  '
Set-TextCell 4 2 'TetherUSD'
Set-TextCell 4 3 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextCell 4 4 '1.003'
Set-TextCell 4 5 '  -0.48%  '
Set-TextCell 5 2 'BNB'
Set-TextCell 5 3 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextCell 5 4 '324.73'
Set-TextCell 5 5 '  -0.91%  '
Set-TextCell 6 2 'USDC'
Set-TextCell 6 3 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextCell 6 4 '1.001'
Set-TextCell 6 5 '  -0.49%  '
Set-TextCell 7 2 'XRP'
Set-TextCell 7 3 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextCell 7 4 '0.4592'
Set-TextCell 7 5 '  -1.72%  '
Set-TextCell 8 2 'Cardano'
Set-TextCell 8 3 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell 8 4 '0.3826'
Set-TextCell 8 5 '  -2.42%  '
Set-TextCell 9 2 'Dogecoin'
Set-TextCell 9 3 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 9 4 '0.07714'
Set-TextCell 9 5 '  -3.05%  '
Set-TextCell 10 2 'Polygon'
Set-TextCell 10 3 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 10 4 '0.9799'
Set-TextCell 10 5 '  -0.99%  '
Set-TextCell 11 2 'Solana'
Set-TextCell 11 3 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell 11 4 '22.07'
Set-TextCell 11 5 '  -3.20%  '
Set-TextCell 12 2 'WrappedEther'
Set-TextCell 12 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 12 4 '1.899.42'
Set-TextCell 12 5 '  -5.44%  '
Set-TextCell 13 2 'Polkadot'
Set-TextCell 13 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 13 4 '5.671'
Set-TextCell 13 5 '  -2.55%  '
Set-TextCell 14 2 'Chainlink'
Set-TextCell 14 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 14 4 '6.933'
Set-TextCell 14 5 '  -3.61%  '
Set-TextCell 15 2 'TRON'
Set-TextCell 15 3 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 15 4 '0.07047'
Set-TextCell 15 5 '  -1.25%  '
Set-TextCell 16 2 'BinanceUSD'
Set-TextCell 16 3 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 16 4 '1.004'
Set-TextCell 16 5 '  -0.29%  '
Set-TextCell 17 2 'Litecoin'
Set-TextCell 17 3 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 17 4 '83.80'
Set-TextCell 17 5 '  -4.71%  '
Set-TextCell 18 2 'ShibaInu'
Set-TextCell 18 3 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 18 4 '0.000009470'
Set-TextCell 18 5 '  -4.72%  '
Set-TextCell 19 2 'Avalanche'
Set-TextCell 19 3 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell 19 4 '16.70'
Set-TextCell 19 5 '  -3.30%  '
Set-TextCell 20 2 'Dai'
Set-TextCell 20 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 20 4 '1.002'
Set-TextCell 20 5 '  -0.41%  '
Set-TextCell 21 2 'WrappedBTC'
Set-TextCell 21 3 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell 21 4 '28.965.65'
Set-TextCell 21 5 '  -2.09%  '
Set-TextCell 22 2 'Uniswap'
Set-TextCell 22 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 22 4 '5.321'
Set-TextCell 22 5 '  -4.19%  '
Set-TextCell 23 2 'Cosmos'
Set-TextCell 23 3 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 23 4 '10.88'
Set-TextCell 23 5 '  -2.56%  '
Set-TextCell 24 2 'WrappedliquidstakedEther2.0'
Set-TextCell 24 3 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 24 4 '2.133.73'
Set-TextCell 24 5 '  -5.28%  '
Set-TextCell 25 2 'Toncoin'
Set-TextCell 25 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 25 4 '2.094'
Set-TextCell 25 5 '  -1.11%  '
Set-TextCell 26 2 'Monero'
Set-TextCell 26 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 26 4 '158.34'
Set-TextCell 26 5 '  -0.41%  '
Set-TextCell 27 2 'EthereumClassic'
Set-TextCell 27 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 27 4 '19.09'
Set-TextCell 27 5 '  -2.46%  '
Set-TextCell 28 2 'InternetComputer(DFINITY)'
Set-TextCell 28 3 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 28 4 '5.660'
Set-TextCell 28 5 '  -2.87%  '
Set-TextCell 29 2 'BitcoinCash'
Set-TextCell 29 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell 29 4 '117.34'
Set-TextCell 29 5 '  -2.27%  '
Set-TextCell 30 2 'LidoDAOToken'
Set-TextCell 30 3 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 30 4 '1.854'
Set-TextCell 30 5 '  -2.30%  '
Set-TextCell 31 2 'Stellar'
Set-TextCell 31 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 31 4 '0.09284'
Set-TextCell 31 5 '  -1.56%  '
Set-TextCell 32 2 'ImmutableX'
Set-TextCell 32 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 32 4 '0.8641'
Set-TextCell 32 5 '  -2.41%  '
Set-TextCell 33 2 'Filecoin'
Set-TextCell 33 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 33 4 '5.072'
Set-TextCell 33 5 '  -3.28%  '
Set-TextCell 34 2 'ARBITRUM'
Set-TextCell 34 3 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 34 4 '1.245'
Set-TextCell 34 5 '  -5.65%  '
Set-TextCell 35 2 'HuobiToken'
Set-TextCell 35 3 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 35 4 '3.017'
Set-TextCell 35 5 '  -4.69%  '
Set-TextCell 36 2 'Hedera'
Set-TextCell 36 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 36 4 '0.05729'
Set-TextCell 36 5 '  -1.59%  '
Set-TextCell 37 2 'TrustWalletToken'
Set-TextCell 37 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 37 4 '1.154'
Set-TextCell 37 5 '  -1.49%  '
Set-TextCell 38 2 'Frax'
Set-TextCell 38 3 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextCell 38 4 '1.002'
Set-TextCell 38 5 '  -0.41%  '
Set-TextCell 39 2 'VeChain'
Set-TextCell 39 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 39 4 '0.02038'
Set-TextCell 39 5 '  -3.30%  '
Set-TextCell 40 2 'TheSandbox'
Set-TextCell 40 3 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 40 4 '0.5508'
Set-TextCell 40 5 '  -3.69%  '
Set-TextCell 41 2 'FraxShare'
Set-TextCell 41 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 41 4 '7.401'
Set-TextCell 41 5 '  -4.41%  '
Set-TextCell 42 2 'Algorand'
Set-TextCell 42 3 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 42 4 '0.1755'
Set-TextCell 42 5 '  -2.41%  '
Set-TextCell 43 2 'MXToken'
Set-TextCell 43 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 43 4 '2.847'
Set-TextCell 43 5 '  +3.41%  '
Set-TextCell 44 2 'Aptos'
Set-TextCell 44 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 44 4 '9.338'
Set-TextCell 44 5 '  -3.26%  '
Set-TextCell 45 2 'Decentraland'
Set-TextCell 45 3 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell 45 4 '0.5185'
Set-TextCell 45 5 '  -2.98%  '
Set-TextCell 46 2 'EnergySwap'
Set-TextCell 46 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 46 4 '11.25'
Set-TextCell 46 5 '  -4.64%  '
Set-TextCell 47 2 'Cronos'
Set-TextCell 47 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 47 4 '0.06828'
Set-TextCell 47 5 '  -1.63%  '
Set-TextCell 48 2 'RenderToken'
Set-TextCell 48 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 48 4 '2.046'
Set-TextCell 48 5 '  -5.01%  '
Set-TextCell 49 2 'Quant'
Set-TextCell 49 3 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 49 4 '110.88'
Set-TextCell 49 5 '  -2.51%  '
Set-TextCell 50 2 'NEARProtocol'
Set-TextCell 50 3 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 50 4 '1.779'
Set-TextCell 50 5 '  -2.91%  '
Set-TextCell 51 2 'PEPE'
Set-TextCell 51 3 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell 51 4 '0.000002557'
Set-TextCell 51 5 '  -5.43%  '
